$wb = $excel.ActiveWorkbook

# --- DATA_SET sheet: remove the rows that correspond to the deleted
#     dosing-error records. Delete from the bottom up so row numbers of
#     rows not yet processed stay stable.
$wsData = $wb.Worksheets.Item("DATA_SET")
$wsData.Rows.Item(8).Delete()
$wsData.Rows.Item(4).Delete()
$wsData.Rows.Item(3).Delete()

# --- PATIENT sheet: fix the creatinine clearance value used by the error model
$wsPatient = $wb.Worksheets.Item("PATIENT")
$wsPatient.Range("C2").Value = 130
